$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"

# --- Cells changing data type (numeric <-> text placeholder) ---
$ws.Range("F16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 1

$ws.Range("F16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 2

$ws.Range("H16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = -50

$ws.Range("C23").Value = "'0"
$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("D27").Value = "'0"
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2

$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 50

# --- Remaining numeric value updates ---
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 55
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = 52.777777777777
$ws.Range("L16").Value = 19.565217391304
$ws.Range("M16").Value = -39.560439560439
$ws.Range("N16").Value = -79.007633587786
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 128.571428571429
$ws.Range("I17").Value = 58
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = -3.333333333333
$ws.Range("L17").Value = -12.121212121212
$ws.Range("M17").Value = 20.833333333333
$ws.Range("N17").Value = -52.066115702479
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = -38.095238095238
$ws.Range("L18").Value = -39.53488372093
$ws.Range("M18").Value = -74.757281553398
$ws.Range("N18").Value = -94.990366088632
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = -2
$ws.Range("I19").Value = 239
$ws.Range("J19").Value = 247
$ws.Range("K19").Value = -3.238866396761
$ws.Range("L19").Value = 32.044198895027
$ws.Range("M19").Value = 23.195876288659
$ws.Range("N19").Value = -14.028776978417
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 60
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = 46.341463414634
$ws.Range("L20").Value = 36.363636363636
$ws.Range("M20").Value = -13.043478260869
$ws.Range("N20").Value = -94.858611825192
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = 18.518518518518
$ws.Range("I21").Value = 441
$ws.Range("J21").Value = 432
$ws.Range("K21").Value = 2.083333333333
$ws.Range("L21").Value = 14.84375
$ws.Range("M21").Value = -13.529411764705
$ws.Range("N21").Value = -81.353065539112
$ws.Range("F23").Value = 3
$ws.Range("L23").Value = 23.076923076923
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 18
$ws.Range("F24").Value = 123
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = 30.851063829787
$ws.Range("I24").Value = 626
$ws.Range("J24").Value = 479
$ws.Range("K24").Value = 30.688935281837
$ws.Range("L24").Value = 48.693586698337
$ws.Range("M24").Value = 69.647696476964
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 73.333333333333
$ws.Range("F25").Value = 100
$ws.Range("G25").Value = 62
$ws.Range("H25").Value = 61.290322580645
$ws.Range("I25").Value = 519
$ws.Range("J25").Value = 315
$ws.Range("K25").Value = 64.761904761904
$ws.Range("L25").Value = 93.65671641791
$ws.Range("C26").Value = 8
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 28
$ws.Range("H26").Value = 55.555555555555
$ws.Range("I26").Value = 122
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 22
$ws.Range("L26").Value = 46.987951807228
$ws.Range("M26").Value = -5.426356589147
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = -72.727272727272
$ws.Range("L27").Value = -62.5
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 14
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = -12.5
$ws.Range("L28").Value = 16.666666666666
$ws.Range("N29").Value = -90.90909090909
$ws.Range("N30").Value = -90
